$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the date in A1 (price list date)
$ws.Range("A1").Value = 45436

# Update price column D for the updated price list
$ws.Range("D19").Value = 3113.977
$ws.Range("D20").Value = 3297.575
$ws.Range("D21").Value = 3702.292
$ws.Range("D22").Value = 4106.981
$ws.Range("D23").Value = 4594.128
$ws.Range("D24").Value = 4961.364
$ws.Range("D25").Value = 5493.473
$ws.Range("D26").Value = 5830.724
$ws.Range("D27").Value = 6445.273
$ws.Range("D28").Value = 7239.687
$ws.Range("D29").Value = 8019.126
$ws.Range("D30").Value = 9068.355
$ws.Range("D31").Value = 10642.201
$ws.Range("D38").Value = 10956.969
$ws.Range("D39").Value = 12920.532
$ws.Range("D40").Value = 15063.945
$ws.Range("D41").Value = 19043.537
$ws.Range("D42").Value = 24484.54
$ws.Range("D43").Value = 30142.901
$ws.Range("D44").Value = 34100.001
$ws.Range("D45").Value = 39196.257
$ws.Range("D46").Value = 17222.362
$ws.Range("D53").Value = 5388.545
$ws.Range("D54").Value = 6475.246
$ws.Range("D55").Value = 6947.398
$ws.Range("D56").Value = 7884.227
$ws.Range("D57").Value = 8575.829
$ws.Range("D58").Value = 9255.710999999999
$ws.Range("D59").Value = 10035.144
$ws.Range("D60").Value = 10349.909
$ws.Range("D61").Value = 11301.703
$ws.Range("D62").Value = 13400.184
$ws.Range("D63").Value = 14464.404
$ws.Range("D64").Value = 17836.919
